# Auto-applied update to match scheduled-runner market data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 14502.571
$ws.Range("I69").Value = 9703.799999999999
$ws.Range("K69").Value = 29111.4
$ws.Range("M69").Value = -28237.4
$ws.Range("H70").Value = 6973.0435
$ws.Range("I70").Value = 2055.5
$ws.Range("J70").Value = 18213.143
$ws.Range("K70").Value = 6166.5
$ws.Range("L70").Value = 54639.429
$ws.Range("M70").Value = -5896.5
$ws.Range("N70").Value = -55179.429
$ws.Range("H72").Value = 14502.571
$ws.Range("I72").Value = 9703.799999999999
$ws.Range("K72").Value = 87334.2
$ws.Range("M72").Value = -82966.2
$ws.Range("H73").Value = 6973.0435
$ws.Range("I73").Value = 2055.5
$ws.Range("J73").Value = 18213.143
$ws.Range("K73").Value = 6166.5
$ws.Range("L73").Value = 54639.429
$ws.Range("M73").Value = -5230.5
$ws.Range("N73").Value = -56511.429
$ws.Range("H113").Value = 3484.35
$ws.Range("J113").Value = 3616
$ws.Range("L113").Value = 3616
$ws.Range("N113").Value = -10124
$ws.Range("H127").Value = 1533.7858
$ws.Range("I127").Value = 1533.7858
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 4601.357400000001
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 358.6425999999992
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 7930.972
$ws.Range("I132").Value = 3057.862
$ws.Range("J132").Value = 28119.572
$ws.Range("K132").Value = 9173.585999999999
$ws.Range("L132").Value = 84358.716
$ws.Range("M132").Value = -6643.585999999999
$ws.Range("N132").Value = -89418.716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5783.5625
$ws.Range("I61").Value = 5795.5
$ws.Range("K61").Value = 5795.5
$ws.Range("M61").Value = -5583.5
$ws.Range("H74").Value = 3398.4167
$ws.Range("I74").Value = 3035.625
$ws.Range("J74").Value = 4124
$ws.Range("K74").Value = 3035.625
$ws.Range("L74").Value = 4124
$ws.Range("M74").Value = -2161.625
$ws.Range("N74").Value = -5872
$ws.Range("H77").Value = 3398.4167
$ws.Range("I77").Value = 3035.625
$ws.Range("J77").Value = 4124
$ws.Range("K77").Value = 15178.125
$ws.Range("L77").Value = 20620
$ws.Range("M77").Value = -10810.125
$ws.Range("N77").Value = -29356
$ws.Range("H136").Value = 5783.5625
$ws.Range("I136").Value = 5795.5
$ws.Range("K136").Value = 17386.5
$ws.Range("M136").Value = -14836.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1855.5
$ws.Range("J80").Value = 1843.1
$ws.Range("L80").Value = 1843.1
$ws.Range("N80").Value = -3839.1
$ws.Range("H83").Value = 1855.5
$ws.Range("J83").Value = 1843.1
$ws.Range("L83").Value = 9215.5
$ws.Range("N83").Value = -19199.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1315
$ws.Range("I22").Value = 385.66666
$ws.Range("K22").Value = 385.66666
$ws.Range("M22").Value = -35.66665999999998
$ws.Range("H31").Value = 2138.7144
$ws.Range("I31").Value = 1932.8889
$ws.Range("J31").Value = 2509.2
$ws.Range("K31").Value = 1932.8889
$ws.Range("L31").Value = 2509.2
$ws.Range("M31").Value = -1637.8889
$ws.Range("N31").Value = -3099.2
$ws.Range("H34").Value = 2138.7144
$ws.Range("I34").Value = 1932.8889
$ws.Range("J34").Value = 2509.2
$ws.Range("K34").Value = 1932.8889
$ws.Range("L34").Value = 2509.2
$ws.Range("M34").Value = -1730.8889
$ws.Range("N34").Value = -2913.2
$ws.Range("H58").Value = 2679.3845
$ws.Range("I58").Value = 2183.3
$ws.Range("K58").Value = 2183.3
$ws.Range("M58").Value = -1980.3
$ws.Range("H99").Value = 10288.182
$ws.Range("J99").Value = 12888
$ws.Range("L99").Value = 12888
$ws.Range("N99").Value = -15884
$ws.Range("H126").Value = 10288.182
$ws.Range("J126").Value = 12888
$ws.Range("L126").Value = 38664
$ws.Range("N126").Value = -43604
$ws.Range("H132").Value = 9911.588
$ws.Range("I132").Value = 7285.7144
$ws.Range("J132").Value = 11749.7
$ws.Range("K132").Value = 21857.1432
$ws.Range("L132").Value = 35249.10000000001
$ws.Range("M132").Value = -19327.1432
$ws.Range("N132").Value = -40309.10000000001
$ws.Range("H134").Value = 3299.1064
$ws.Range("I134").Value = 3181.4866
$ws.Range("J134").Value = 3734.3
$ws.Range("K134").Value = 9544.459800000001
$ws.Range("L134").Value = 11202.9
$ws.Range("M134").Value = -7009.459800000001
$ws.Range("N134").Value = -16272.9
$ws.Range("H136").Value = 2679.3845
$ws.Range("I136").Value = 2183.3
$ws.Range("K136").Value = 6549.900000000001
$ws.Range("M136").Value = -3999.900000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 228515.69
$ws.Range("J46").Value = 357977.34
$ws.Range("L46").Value = 1073932.02
$ws.Range("N46").Value = -1074114.02
$ws.Range("H125").Value = 8288.333000000001
$ws.Range("J125").Value = 9933
$ws.Range("L125").Value = 29799
$ws.Range("N125").Value = -39639

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4294.0435
$ws.Range("I132").Value = 4335.5
$ws.Range("K132").Value = 13006.5
$ws.Range("M132").Value = -10476.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2808.5
$ws.Range("I7").Value = 2384.5715
$ws.Range("J7").Value = 3797.6667
$ws.Range("K7").Value = 2384.5715
$ws.Range("L7").Value = 3797.6667
$ws.Range("M7").Value = -2272.5715
$ws.Range("N7").Value = -4021.6667
$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("N44").Value = -50912
$ws.Range("H126").Value = 2808.5
$ws.Range("I126").Value = 2384.5715
$ws.Range("J126").Value = 3797.6667
$ws.Range("K126").Value = 7153.7145
$ws.Range("L126").Value = 11393.0001
$ws.Range("M126").Value = -4683.7145
$ws.Range("N126").Value = -16333.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 35970
$ws.Range("J54").Value = 45606.668
$ws.Range("L54").Value = 45606.668
$ws.Range("N54").Value = -46646.668
$ws.Range("H61").Value = 8449
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 8449
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 8449
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -9033
$ws.Range("H62").Value = 2998.5
$ws.Range("I62").Value = 2998.5
$ws.Range("K62").Value = 2998.5
$ws.Range("M62").Value = -2374.5
$ws.Range("H65").Value = 2998.5
$ws.Range("I65").Value = 2998.5
$ws.Range("K65").Value = 14992.5
$ws.Range("M65").Value = -11872.5
$ws.Range("H70").Value = 48998
$ws.Range("I70").Value = 48998
$ws.Range("K70").Value = 48998
$ws.Range("M70").Value = -48683
$ws.Range("H73").Value = 48998
$ws.Range("I73").Value = 48998
$ws.Range("K73").Value = 48998
$ws.Range("M73").Value = -47906
$ws.Range("H107").Value = 596.9474
$ws.Range("I107").Value = 245.84616
$ws.Range("J107").Value = 1357.6666
$ws.Range("K107").Value = 737.5384799999999
$ws.Range("L107").Value = 4072.9998
$ws.Range("M107").Value = 1182.46152
$ws.Range("N107").Value = -7912.9998
$ws.Range("H113").Value = 3197.5881
$ws.Range("I113").Value = 2759.3845
$ws.Range("K113").Value = 8278.1535
$ws.Range("M113").Value = -6108.1535
$ws.Range("H122").Value = 6407.9443
$ws.Range("I122").Value = 2577
$ws.Range("K122").Value = 7731
$ws.Range("M122").Value = -5281
$ws.Range("H132").Value = 3135.8718
$ws.Range("I132").Value = 2818.2188
$ws.Range("K132").Value = 8454.6564
$ws.Range("M132").Value = -5924.6564
$ws.Range("H136").Value = 1554.9412
$ws.Range("I136").Value = 1323.6511
$ws.Range("K136").Value = 3970.9533
$ws.Range("M136").Value = -1420.9533
